$d = $word.ActiveDocument

# --- Edit 1: extend the first paragraph with a red "(This is a change ... )" note ---
# Original run: "This is a Microsoft word document." (positions 0-34, paragraph mark at 35)
$p1End = $d.Paragraphs(1).Range.End - 1

# Append two trailing spaces to the existing (uncolored) run.
$rSpaces = $d.Range($p1End, $p1End)
$rSpaces.InsertAfter("  ")

# Run 2 (red): "(This is a change " + en dash + " Ve"
$pos1 = $rSpaces.End
$rNote1 = $d.Range($pos1, $pos1)
$rNote1.InsertAfter("(This is a change – Ve")
$rNote1.Font.Color = 192

# Run 3 (red): "rsion for branch alternate"
$pos2 = $rNote1.End
$rNote2 = $d.Range($pos2, $pos2)
$rNote2.InsertAfter("rsion for branch alternate")
$rNote2.Font.Color = 192

# Run 4 (red): ")"
$pos3 = $rNote2.End
$rNote3 = $d.Range($pos3, $pos3)
$rNote3.InsertAfter(")")
$rNote3.Font.Color = 192

# --- Edit 2: append a new, blank, shaded paragraph after the last paragraph ---
$tailRange = $d.Content
$tailRange.Find.Execute(
    "we are free at last.", $false, $false, $false, $false, $false,
    $true, 1, $false, "we are free at last.^p", 2)

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Style = "Normal"
$newPara.Shading.Texture = 0
$newPara.Shading.ForegroundPatternColor = -16777216
$newPara.Shading.BackgroundPatternColor = 16382457
